# Insert a new data row (row 31) into the sheet, shifting the existing
# rows 31-70 down to 32-71. This grows the used range from A1:R70 to
# A1:R71 and adds one new record ("$/saco 25 kilos" / Región del Maule)
# for Ají - Americana (o) while leaving every other existing row intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 31..70 down to 32..71, creating a blank row 31.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with its data.
$ws.Cells.Item(31, 1).Value  = 11
$ws.Cells.Item(31, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(31, 3).Value  = "Bíobío"
$ws.Cells.Item(31, 4).Value  = 44601
$ws.Cells.Item(31, 5).Value  = 8
$ws.Cells.Item(31, 6).Value  = 100112021
$ws.Cells.Item(31, 7).Value  = "Ají"
$ws.Cells.Item(31, 8).Value  = "Americana (o)"
$ws.Cells.Item(31, 9).Value  = "Primera"
$ws.Cells.Item(31, 10).Value = 30
$ws.Cells.Item(31, 11).Value = 18000
$ws.Cells.Item(31, 12).Value = 19000
$ws.Cells.Item(31, 13).Value = 18333
$ws.Cells.Item(31, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Región del Maule"
$ws.Cells.Item(31, 16).Value = 733
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
